$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Unmerge every merged range on the sheet (B1:C1, D1:E1, F1:G1, H1:I1, J1:K1, L1:M1, N1:O1)
$ws.Cells.UnMerge()

# 2. Drop the centered-alignment style that used to be applied to the header
#    cells B1:O1 so they fall back to the default (unstyled) xf.
$ws.Range("B1:O1").ClearFormats()

# 3. Rewrite row 1 (header) - each "Партия N" header is now split into two
#    columns: "Партия N 1" and "Партия N 2" (one exception: "Партия 22"
#    instead of "Партия 2 2").
$ws.Range("A1").Value = "Имя 1"
$ws.Range("B1").Value = "Партия 1 1"
$ws.Range("C1").Value = "Партия 1 2"
$ws.Range("D1").Value = "Партия 2 1"
$ws.Range("E1").Value = "Партия 22"
$ws.Range("F1").Value = "Партия 3 1"
$ws.Range("G1").Value = "Партия 3 2"
$ws.Range("H1").Value = "Партия 4 1"
$ws.Range("I1").Value = "Партия 4 2"
$ws.Range("J1").Value = "Партия 5 1"
$ws.Range("K1").Value = "Партия 5 2"
$ws.Range("L1").Value = "Партия 6 1"
$ws.Range("M1").Value = "Партия 6 2"
$ws.Range("N1").Value = "Партия 7 1"
$ws.Range("O1").Value = "Партия 7 2"
$ws.Range("P1").Value = "Имя 2"
$ws.Range("Q1").Value = "Общий счет"

# Row 2 data values are unchanged; re-assert them for safety/no-op parity.
$ws.Range("A2").Value = "Петя"
$ws.Range("B2").Value = 11
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 12
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = 13
$ws.Range("G2").Value = 11
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "Витя"
$ws.Range("Q2").Value = "'3:1"

# 4. Column widths: A and P narrow, B:O wider (matches the new two-column
#    layout), Q keeps its existing width. (5.5 / 9.5 are the COM
#    "characters" inputs that land closest to the author's target widths
#    of 6.28515625 / 10.28515625 once Excel quantizes them.)
$ws.Columns.Item(1).ColumnWidth = 5.5
$ws.Range("B1:O1").ColumnWidth = 9.5
$ws.Columns.Item(16).ColumnWidth = 5.5

# 5. Move the active selection to M7, as recorded in the saved view state.
$ws.Range("M7").Select() | Out-Null
